$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.999999900753416,
    0.6876245000598846,
    0.999999999999147,
    0.9999999998272675,
    0.9999999998963385,
    0.00000005891705945054375,
    0.1854395906834327,
    0.0000000000001087298210039851,
    0.0000000001042355012398326,
    0.00000000005217211553041828,
    0.00009931498004524596,
    0.0002427283655664161,
    1.000000140112824,
    0.0002530618146055635,
    115.2942703081797,
    165.268179127776
)

# Columns B (2) through Q (17)
for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}
